$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting rows 33:72 down to 34:73
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new record
$ws.Cells.Item(33, 1).Value = 10
$ws.Cells.Item(33, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value = "La Araucanía"
$ws.Cells.Item(33, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = 300000001
$ws.Cells.Item(33, 7).Value = "Rabanito"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 50
$ws.Cells.Item(33, 11).Value = 8000
$ws.Cells.Item(33, 12).Value = 8000
$ws.Cells.Item(33, 13).Value = 8000
$ws.Cells.Item(33, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(33, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(33, 16).Value = 667
$ws.Cells.Item(33, 17).Value = 12
$ws.Cells.Item(33, 18).Value = "Hortaliza"
